$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103; this shifts existing rows 103-205 down to 104-206
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new record's data
$ws.Cells.Item(103, 1).Value = 11
$ws.Cells.Item(103, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(103, 3).Value = "Bíobío"
$ws.Cells.Item(103, 4).Value = (Get-Date -Year 2023 -Month 5 -Day 3 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(103, 5).Value = 8
$ws.Cells.Item(103, 6).Value = 100112043
$ws.Cells.Item(103, 7).Value = "Pepino ensalada"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 100
$ws.Cells.Item(103, 11).Value = 11000
$ws.Cells.Item(103, 12).Value = 12000
$ws.Cells.Item(103, 13).Value = 11500
$ws.Cells.Item(103, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(103, 15).Value = "Región Metropolitana"
$ws.Cells.Item(103, 16).Value = 192
$ws.Cells.Item(103, 17).Value = 60
$ws.Cells.Item(103, 18).Value = "Hortaliza"
